$wb = $excel.ActiveWorkbook

# The "Sage.X3.ReservedSheet" tab (the dictionary/placeholder sheet) gets a
# new placeholder row enabling locale support in non-template mode.
$ws = $wb.Worksheets.Item("Sage.X3.ReservedSheet")

# New placeholder row: A6 stays empty, A7 carries the new token.
$ws.Range("A7").Value = "___addSupportedLocales___"

# Widen column A so the new (longer) placeholder text fits.
$ws.Columns.Item(1).ColumnWidth = 30

# Make this sheet the active one, with B7 selected, matching the user's
# last action after typing the new row.
$ws.Activate()
$ws.Range("B7").Select()
